# Rotate the "Recorded By" (column G) values: move the first
# comma-separated entry to the end of the list, for every data row
# on the active worksheet -- except the specific literal value
# "System, admin@admin.com" which must stay untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

# Column G is "Recorded By" (row 1 is the header row).
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $value = $cell.Value2

    if ($null -eq $value) { continue }
    if ($value -eq "") { continue }

    # Leave this exact value alone.
    if ($value -eq "System, admin@admin.com") { continue }

    if ($value -notmatch ",") { continue }

    $parts = $value -split ",\s*"
    if ($parts.Count -lt 2) { continue }

    $rotated = ($parts[1..($parts.Count - 1)] + $parts[0]) -join ", "

    $cell.Value2 = $rotated
}
